# Team Project Diary - log the latest hours worked, remove reliance on the
# old log and bring the new Excel log up to date across several weeks.
#
# The order the cells below are written in matters: it reproduces the
# order in which new entries were typed into the shared-strings table.
$wb = $excel.ActiveWorkbook

$week1 = $wb.Worksheets.Item("Week 1")
$week2 = $wb.Worksheets.Item("Week 2")
$week3 = $wb.Worksheets.Item("Week 3")
$week4 = $wb.Worksheets.Item("Week 4")
$week5 = $wb.Worksheets.Item("Week 5")

# Week 2 - intro for the proposal
$week2.Range("B27").Value = "Wrote the intro for Proposal"
$week2.Range("C27").Value = 2

# Week 1 - brainstorm meeting
$week1.Range("B26").Value = "Project Brainstorm Meeting"
$week1.Range("C26").Value = 1.5

# Week 4 - updating the log book
$week4.Range("B27").Value = "Updating log book"
$week4.Range("C27").Value = 1

# Week 4 - survey question creation meeting
$week4.Range("B26").Value = "Survey question creation meeting (Telecommute from class)"
$week4.Range("C26").Value = 2

# Week 4 - Friday meeting to finish survey questions
$week4.Range("B28").Value = "Friday meeting to complete survey questions and Milestone confirmation"
$week4.Range("C28").Value = 2

# Week 5 - improving the proposal based on feedback
$week5.Range("B24").Value = "Improving my portion of the proposal based on feedback"
$week5.Range("C24").Value = 1.5

# Week 2 - reading articles
$week2.Range("B26").Value = "Read Articles"
$week2.Range("C26").Value = 2

# Week 3 - presentation meeting
$week3.Range("B24").Value = "Presentation Meeting (Tuesday)"
$week3.Range("C24").Value = 3

# Update each sheet's on-screen selection to match where the author left off
$week1.Activate()
$week1.Range("C26").Select() | Out-Null

$week2.Activate()
$week2.Range("B28").Select() | Out-Null

$week3.Activate()
$week3.Range("B32").Select() | Out-Null

$week4.Activate()
$week4.Range("C29").Select() | Out-Null

# Week 5 becomes the active sheet (author finished here, replacing the old
# standalone log file with this updated workbook)
$week5.Activate()
$week5.Range("B24").Select() | Out-Null
